$wb = $excel.ActiveWorkbook

# --- Sheet "Logs": append new row 12 with testmail #7 data ---
$logs = $wb.Worksheets.Item("Logs")

$logs.Cells.Item(12, 1).Value = "Is dit artikel nog op voorraad?"
$logs.Cells.Item(12, 2).Value = "mailmind.test@zohomail.eu"
$logs.Cells.Item(12, 3).Value = "Testmail #7: Is dit artikel nog op voorraad?"
$logs.Cells.Item(12, 4).Value = "Productinformatie"
$logs.Cells.Item(12, 5).Value = "Bedankt, we hebben dit doorgestuurd naar support@bedrijf.nl."
$logs.Cells.Item(12, 6).Value = "2025-08-01 23:49:37"
$logs.Cells.Item(12, 7).Value = "Ja"
$logs.Cells.Item(12, 8).Value = "Ja"
$logs.Cells.Item(12, 9).Value = "Nee"
$logs.Cells.Item(12, 10).Value = "Nee"

# Extend the conditional formatting ranges (D,G,H,I,J) from row 11 to row 12
foreach ($col in @("D","G","H","I","J")) {
    $oldRange = $col + "2:" + $col + "11"
    $newRange = $col + "2:" + $col + "12"
    $fcs = $logs.Range($oldRange).FormatConditions
    for ($i = 1; $i -le $fcs.Count; $i++) {
        $fc = $fcs.Item($i)
        $fc.ModifyAppliesToRange($logs.Range($newRange))
    }
}

# --- Sheet "Dashboard": swap category rows 4 and 5 ---
$dash = $wb.Worksheets.Item("Dashboard")

$dash.Cells.Item(4, 1).Value = "Productinformatie"
$dash.Cells.Item(4, 2).Value = 2
$dash.Cells.Item(5, 1).Value = "Bestelling / Levering"
$dash.Cells.Item(5, 2).Value = 1
